$wb = $excel.ActiveWorkbook

# Mapping of row -> new F-column value, applied identically to sheets "展览" and "全部类型"
$updates = @{
    3  = 293
    4  = 11135
    5  = 10352
    8  = 728
    10 = 18
    13 = 9619
    14 = 2216
    17 = 39
    22 = 10789
    24 = 17
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
